$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 87.28570999999999
$ws.Range("I6").Value = 68.5
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 205.5
$ws.Range("L6").Value = 600
$ws.Range("M6").Value = -93.5
$ws.Range("N6").Value = -824
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("H15").Value = 1723.3226
$ws.Range("I15").Value = 1723.3226
$ws.Range("K15").Value = 5169.9678
$ws.Range("M15").Value = -5000.9678
$ws.Range("H39").Value = 2654.0908
$ws.Range("I39").Value = 3198.1428
$ws.Range("J39").Value = 1702
$ws.Range("K39").Value = 9594.428400000001
$ws.Range("L39").Value = 5106
$ws.Range("M39").Value = -9298.428400000001
$ws.Range("N39").Value = -5698
$ws.Range("H43").Value = 6040.5
$ws.Range("J43").Value = 6455.25
$ws.Range("L43").Value = 6455.25
$ws.Range("N43").Value = -6593.25
$ws.Range("H82").Value = 6363.25
$ws.Range("I82").Value = 6363.25
$ws.Range("K82").Value = 19089.75
$ws.Range("M82").Value = -18683.75
$ws.Range("H85").Value = 6363.25
$ws.Range("I85").Value = 6363.25
$ws.Range("K85").Value = 19089.75
$ws.Range("M85").Value = -17685.75
$ws.Range("H94").Value = 2549.6667
$ws.Range("I94").Value = 3199.5
$ws.Range("K94").Value = 3199.5
$ws.Range("M94").Value = -2748.5
$ws.Range("H115").Value = 3389.8333
$ws.Range("I115").Value = 1585.5
$ws.Range("J115").Value = 6998.5
$ws.Range("K115").Value = 4756.5
$ws.Range("L115").Value = 20995.5
$ws.Range("M115").Value = -3189.5
$ws.Range("N115").Value = -24129.5
$ws.Range("H132").Value = 316831.94
$ws.Range("I132").Value = 391164.2
$ws.Range("J132").Value = 10211.5
$ws.Range("K132").Value = 1173492.6
$ws.Range("L132").Value = 30634.5
$ws.Range("M132").Value = -1170962.6
$ws.Range("N132").Value = -35694.5
$ws.Range("H137").Value = 10908.421
$ws.Range("I137").Value = 7450.5835
$ws.Range("J137").Value = 16836.143
$ws.Range("K137").Value = 22351.7505
$ws.Range("L137").Value = 50508.429
$ws.Range("M137").Value = -19801.7505
$ws.Range("N137").Value = -55608.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 376.2
$ws.Range("I5").Value = 410.66666
$ws.Range("J5").Value = 324.5
$ws.Range("K5").Value = 410.66666
$ws.Range("L5").Value = 324.5
$ws.Range("M5").Value = -298.66666
$ws.Range("N5").Value = -548.5
$ws.Range("H32").Value = 1524289.4
$ws.Range("I32").Value = 1017.5
$ws.Range("J32").Value = 8554775
$ws.Range("K32").Value = 1017.5
$ws.Range("L32").Value = 8554775
$ws.Range("M32").Value = -730.5
$ws.Range("N32").Value = -8555349
$ws.Range("H61").Value = 5733.5186
$ws.Range("I61").Value = 5533.5415
$ws.Range("K61").Value = 5533.5415
$ws.Range("M61").Value = -5321.5415
$ws.Range("H132").Value = 924568.9
$ws.Range("I132").Value = 1123988.6
$ws.Range("J132").Value = 155378.42
$ws.Range("K132").Value = 3371965.8
$ws.Range("L132").Value = 466135.26
$ws.Range("M132").Value = -3369435.8
$ws.Range("N132").Value = -471195.26
$ws.Range("H136").Value = 5733.5186
$ws.Range("I136").Value = 5533.5415
$ws.Range("K136").Value = 16600.6245
$ws.Range("M136").Value = -14050.6245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 376.2
$ws.Range("I4").Value = 410.66666
$ws.Range("J4").Value = 324.5
$ws.Range("K4").Value = 410.66666
$ws.Range("L4").Value = 324.5
$ws.Range("M4").Value = -295.66666
$ws.Range("N4").Value = -554.5
$ws.Range("H106").Value = 28500
$ws.Range("J106").Value = 28500
$ws.Range("L106").Value = 28500
$ws.Range("N106").Value = -31024
$ws.Range("H134").Value = 1018277.75
$ws.Range("I134").Value = 1165778.2
$ws.Range("J134").Value = 10357.667
$ws.Range("K134").Value = 3497334.6
$ws.Range("L134").Value = 31073.001
$ws.Range("M134").Value = -3494799.6
$ws.Range("N134").Value = -36143.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3914822.5
$ws.Range("J22").Value = 238990
$ws.Range("L22").Value = 238990
$ws.Range("N22").Value = -239690
$ws.Range("H58").Value = 40007972
$ws.Range("I58").Value = 43483556
$ws.Range("K58").Value = 43483556
$ws.Range("M58").Value = -43483353
$ws.Range("H62").Value = 5711.077
$ws.Range("I62").Value = 6538.6665
$ws.Range("J62").Value = 3849
$ws.Range("K62").Value = 6538.6665
$ws.Range("L62").Value = 3849
$ws.Range("M62").Value = -5914.6665
$ws.Range("N62").Value = -5097
$ws.Range("H65").Value = 5711.077
$ws.Range("I65").Value = 6538.6665
$ws.Range("J65").Value = 3849
$ws.Range("K65").Value = 32693.3325
$ws.Range("L65").Value = 19245
$ws.Range("M65").Value = -29573.3325
$ws.Range("N65").Value = -25485
$ws.Range("H69").Value = 60000
$ws.Range("J69").Value = 60000
$ws.Range("L69").Value = 60000
$ws.Range("N69").Value = -61498
$ws.Range("H72").Value = 60000
$ws.Range("J72").Value = 60000
$ws.Range("L72").Value = 180000
$ws.Range("N72").Value = -187488
$ws.Range("H134").Value = 55563670
$ws.Range("I134").Value = 66672200
$ws.Range("K134").Value = 200016600
$ws.Range("M134").Value = -200014065
$ws.Range("H136").Value = 40007972
$ws.Range("I136").Value = 43483556
$ws.Range("K136").Value = 130450668
$ws.Range("M136").Value = -130448118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1222.4286
$ws.Range("I34").Value = 501.83334
$ws.Range("J34").Value = 1762.875
$ws.Range("K34").Value = 1505.50002
$ws.Range("L34").Value = 5288.625
$ws.Range("M34").Value = -1421.50002
$ws.Range("N34").Value = -5456.625
$ws.Range("H52").Value = 1390272.8
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 1390272.8
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 4170818.4
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -4171350.4
$ws.Range("H109").Value = 11709.333
$ws.Range("I109").Value = 9164
$ws.Range("K109").Value = 27492
$ws.Range("M109").Value = -26452
$ws.Range("H134").Value = 14946
$ws.Range("I134").Value = 10231.6
$ws.Range("K134").Value = 30694.8
$ws.Range("M134").Value = -25624.8
$ws.Range("H139").Value = 23811914
$ws.Range("I139").Value = 33335322
$ws.Range("K139").Value = 100005966
$ws.Range("M139").Value = -100000826
$ws.Range("H140").Value = 32830326
$ws.Range("J140").Value = 3698.4546
$ws.Range("L140").Value = 11095.3638
$ws.Range("N140").Value = -21455.3638

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19580.4
$ws.Range("I7").Value = 29018
$ws.Range("J7").Value = 15535.714
$ws.Range("K7").Value = 29018
$ws.Range("L7").Value = 15535.714
$ws.Range("M7").Value = -28906
$ws.Range("N7").Value = -15759.714
$ws.Range("H40").Value = 7004.9546
$ws.Range("I40").Value = 5903.385
$ws.Range("J40").Value = 8596.111000000001
$ws.Range("K40").Value = 5903.385
$ws.Range("L40").Value = 8596.111000000001
$ws.Range("M40").Value = -5767.385
$ws.Range("N40").Value = -8868.111000000001
$ws.Range("H126").Value = 19580.4
$ws.Range("I126").Value = 29018
$ws.Range("J126").Value = 15535.714
$ws.Range("K126").Value = 87054
$ws.Range("L126").Value = 46607.142
$ws.Range("M126").Value = -84584
$ws.Range("N126").Value = -51547.142
$ws.Range("H132").Value = 7004.08
$ws.Range("I132").Value = 4858
$ws.Range("K132").Value = 14574
$ws.Range("M132").Value = -12044
$ws.Range("H136").Value = 62505224
$ws.Range("I136").Value = 26320718
$ws.Range("J136").Value = 200006340
$ws.Range("K136").Value = 78962154
$ws.Range("L136").Value = 600019020
$ws.Range("M136").Value = -78959604
$ws.Range("N136").Value = -600024120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 18416.5
$ws.Range("I62").Value = 21125
$ws.Range("J62").Value = 17062.25
$ws.Range("K62").Value = 21125
$ws.Range("L62").Value = 17062.25
$ws.Range("M62").Value = -20501
$ws.Range("N62").Value = -18310.25
$ws.Range("H65").Value = 18416.5
$ws.Range("I65").Value = 21125
$ws.Range("J65").Value = 17062.25
$ws.Range("K65").Value = 105625
$ws.Range("L65").Value = 85311.25
$ws.Range("M65").Value = -102505
$ws.Range("N65").Value = -91551.25
$ws.Range("H81").Value = 1335.68
$ws.Range("I81").Value = 1226.909
$ws.Range("J81").Value = 2133.3333
$ws.Range("K81").Value = 2453.818
$ws.Range("L81").Value = 4266.6666
$ws.Range("M81").Value = -1392.818
$ws.Range("N81").Value = -6388.6666
$ws.Range("H84").Value = 1335.68
$ws.Range("I84").Value = 1226.909
$ws.Range("J84").Value = 2133.3333
$ws.Range("K84").Value = 12269.09
$ws.Range("L84").Value = 21333.333
$ws.Range("M84").Value = -6965.09
$ws.Range("N84").Value = -31941.333
$ws.Range("H132").Value = 8875.913
$ws.Range("I132").Value = 4954.3335
$ws.Range("J132").Value = 22993.6
$ws.Range("K132").Value = 14863.0005
$ws.Range("L132").Value = 68980.79999999999
$ws.Range("M132").Value = -12333.0005
$ws.Range("N132").Value = -74040.79999999999
$ws.Range("H136").Value = 17248126
$ws.Range("I136").Value = 23811036
$ws.Range("J136").Value = 20492.25
$ws.Range("K136").Value = 71433108
$ws.Range("L136").Value = 61476.75
$ws.Range("M136").Value = -71430558
$ws.Range("N136").Value = -66576.75
